$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 147, shifting existing rows 147-213 down to 148-214
$ws.Rows.Item(147).EntireRow.Insert()

# Populate the newly inserted row 147 with the new price entry
$ws.Range("A147").Value = 3
$ws.Range("B147").Value = "Femacal de La Calera"
$ws.Range("C147").Value = "Coquimbo"
$ws.Range("D147").Value = 44455
$ws.Range("E147").Value = 5
$ws.Range("F147").Value = 100112003
$ws.Range("G147").Value = "Ajo"
$ws.Range("H147").Value = "Chino"
$ws.Range("I147").Value = "Primera"
$ws.Range("J147").Value = 70
$ws.Range("K147").Value = 15500
$ws.Range("L147").Value = 16000
$ws.Range("M147").Value = 15714
$ws.Range("N147").Value = "$/caja 10 kilos"
$ws.Range("O147").Value = "China"
$ws.Range("P147").Value = 1571
$ws.Range("Q147").Value = 10
$ws.Range("R147").Value = "Hortaliza"
